# Auto-generated cell updates reproducing the Hyperion_Profits scheduled-runner refresh.
# For every affected row we rewrite the calculated price/profit columns (H:N) with
# their freshly recomputed values; a couple of rows gain/lose a trailing column
# (N) because HQ pricing data became available/unavailable for that leve item.
$wb = $excel.ActiveWorkbook

# ========== ALC ==========
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1753
$ws.Range("I19").Value = 921.25
$ws.Range("J19").Value = 2418.4
$ws.Range("K19").Value = 921.25
$ws.Range("L19").Value = 2418.4
$ws.Range("M19").Value = -746.25
$ws.Range("N19").Value = -2768.4
# Row 32
$ws.Range("H32").Value = 5045.107
$ws.Range("J32").Value = 4787.5186
$ws.Range("L32").Value = 4787.5186
$ws.Range("N32").Value = -5439.5186
# Row 41
$ws.Range("H41").Value = 15873855
$ws.Range("I41").Value = 480.33334
$ws.Range("K41").Value = 480.33334
$ws.Range("M41").Value = -40.33334000000002
# Row 53
$ws.Range("H53").Value = 4745.609
$ws.Range("I53").Value = 197.07692
$ws.Range("J53").Value = 10658.7
$ws.Range("K53").Value = 197.07692
$ws.Range("L53").Value = 10658.7
$ws.Range("M53").Value = 439.92308
$ws.Range("N53").Value = -11932.7
# Row 132
$ws.Range("H132").Value = 29415162
$ws.Range("I132").Value = 35718056
$ws.Range("K132").Value = 107154168
$ws.Range("M132").Value = -107151638
# Row 133
$ws.Range("H133").Value = 49999.332
$ws.Range("J133").Value = 49999.332
$ws.Range("L133").Value = 49999.332
$ws.Range("N133").Value = -60119.332
# Row 135
$ws.Range("H135").Value = 845.35
$ws.Range("J135").Value = 1058
$ws.Range("L135").Value = 9522
$ws.Range("N135").Value = -14592
# Row 136
$ws.Range("H136").Value = 195786.12
$ws.Range("J136").Value = 195786.12
$ws.Range("L136").Value = 195786.12
$ws.Range("N136").Value = -205986.12
# Row 137
$ws.Range("H137").Value = 79341.52
$ws.Range("I137").Value = 112387.31
$ws.Range("K137").Value = 337161.93
$ws.Range("M137").Value = -334611.93
# Row 141
$ws.Range("H141").Value = 12613.615
$ws.Range("I141").Value = 6165.737
$ws.Range("K141").Value = 18497.211
$ws.Range("M141").Value = -13317.211

# ========== ARM ==========
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5104.019
$ws.Range("I32").Value = 3836.6099
$ws.Range("J32").Value = 9828
$ws.Range("K32").Value = 3836.6099
$ws.Range("L32").Value = 9828
$ws.Range("M32").Value = -3549.6099
$ws.Range("N32").Value = -10402
# Row 45
$ws.Range("H45").Value = 7575364.5
$ws.Range("I45").Value = 11067418
$ws.Range("J45").Value = 9248.666999999999
$ws.Range("K45").Value = 11067418
$ws.Range("L45").Value = 9248.666999999999
$ws.Range("M45").Value = -11067041
$ws.Range("N45").Value = -10002.667
# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21538
# Row 61
$ws.Range("H61").Value = 6538.773
$ws.Range("I61").Value = 7300.1665
$ws.Range("J61").Value = 3112.5
$ws.Range("K61").Value = 7300.1665
$ws.Range("L61").Value = 3112.5
$ws.Range("M61").Value = -7088.1665
$ws.Range("N61").Value = -3536.5
# Row 74
$ws.Range("H74").Value = 174462.1
$ws.Range("I74").Value = 17536.875
$ws.Range("J74").Value = 592929.3
$ws.Range("K74").Value = 17536.875
$ws.Range("L74").Value = 592929.3
$ws.Range("M74").Value = -16662.875
$ws.Range("N74").Value = -594677.3
# Row 77
$ws.Range("H77").Value = 174462.1
$ws.Range("I77").Value = 17536.875
$ws.Range("J77").Value = 592929.3
$ws.Range("K77").Value = 87684.375
$ws.Range("L77").Value = 2964646.5
$ws.Range("M77").Value = -83316.375
$ws.Range("N77").Value = -2973382.5
# Row 136
$ws.Range("H136").Value = 6538.773
$ws.Range("I136").Value = 7300.1665
$ws.Range("J136").Value = 3112.5
$ws.Range("K136").Value = 21900.4995
$ws.Range("L136").Value = 9337.5
$ws.Range("M136").Value = -19350.4995
$ws.Range("N136").Value = -14437.5

# ========== BSM ==========
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2876.9
$ws.Range("I20").Value = 2703.375
$ws.Range("J20").Value = 2992.5833
$ws.Range("K20").Value = 2703.375
$ws.Range("L20").Value = 2992.5833
$ws.Range("M20").Value = -2456.375
$ws.Range("N20").Value = -3486.5833

# ========== CRP ==========
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3004.8572
$ws.Range("I16").Value = 2572.25
$ws.Range("K16").Value = 2572.25
$ws.Range("M16").Value = -2285.25
# Row 31
$ws.Range("H31").Value = 15510.651
$ws.Range("I31").Value = 8135.4375
$ws.Range("K31").Value = 8135.4375
$ws.Range("M31").Value = -7840.4375
# Row 34
$ws.Range("H34").Value = 15510.651
$ws.Range("I34").Value = 8135.4375
$ws.Range("K34").Value = 8135.4375
$ws.Range("M34").Value = -7933.4375
# Row 105
$ws.Range("H105").Value = 2070.524
$ws.Range("I105").Value = 2070.524
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2070.524
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -323.5239999999999
$ws.Range("N105").ClearContents()
# Row 113
$ws.Range("H113").Value = 3004.8572
$ws.Range("I113").Value = 2572.25
$ws.Range("K113").Value = 2572.25
$ws.Range("M113").Value = -402.25
# Row 122
$ws.Range("H122").Value = 2253.8823
$ws.Range("I122").Value = 1847.9286
$ws.Range("J122").Value = 4148.3335
$ws.Range("K122").Value = 5543.7858
$ws.Range("L122").Value = 12445.0005
$ws.Range("M122").Value = -3093.7858
$ws.Range("N122").Value = -17345.0005
# Row 132
$ws.Range("H132").Value = 74088.57000000001
$ws.Range("I132").Value = 79403.08
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 238209.24
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -235679.24
$ws.Range("N132").Value = -20060

# ========== CUL ==========
$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 125373.5
$ws.Range("I44").Value = 414.66666
$ws.Range("K44").Value = 1243.99998
$ws.Range("M44").Value = -845.9999800000001
# Row 68
$ws.Range("H68").Value = 540
$ws.Range("I68").Value = 551.25
$ws.Range("J68").Value = 495
$ws.Range("K68").Value = 1653.75
$ws.Range("L68").Value = 1485
$ws.Range("M68").Value = -842.75
$ws.Range("N68").Value = -3107
# Row 71
$ws.Range("H71").Value = 540
$ws.Range("I71").Value = 551.25
$ws.Range("J71").Value = 495
$ws.Range("K71").Value = 4961.25
$ws.Range("L71").Value = 4455
$ws.Range("M71").Value = -905.25
$ws.Range("N71").Value = -12567
# Row 122
$ws.Range("H122").Value = 512.0714
$ws.Range("I122").Value = 744.75
$ws.Range("J122").Value = 201.83333
$ws.Range("K122").Value = 6702.75
$ws.Range("L122").Value = 1816.49997
$ws.Range("M122").Value = -4252.75
$ws.Range("N122").Value = -6716.49997
# Row 132
$ws.Range("H132").Value = 1914.9166
$ws.Range("I132").Value = 949
$ws.Range("J132").Value = 2236.889
$ws.Range("K132").Value = 8541
$ws.Range("L132").Value = 20132.001
$ws.Range("M132").Value = -6011
$ws.Range("N132").Value = -25192.001

# ========== WVR ==========
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10106.477
$ws.Range("I62").Value = 13694.1
$ws.Range("K62").Value = 13694.1
$ws.Range("M62").Value = -13070.1
# Row 65
$ws.Range("H65").Value = 10106.477
$ws.Range("I65").Value = 13694.1
$ws.Range("K65").Value = 68470.5
$ws.Range("M65").Value = -65350.5
# Row 122
$ws.Range("H122").Value = 4248.75
$ws.Range("I122").Value = 4248.75
$ws.Range("K122").Value = 12746.25
$ws.Range("M122").Value = -10296.25
# Row 132
$ws.Range("H132").Value = 33708550
$ws.Range("I132").Value = 41677500
$ws.Range("K132").Value = 125032500
$ws.Range("M132").Value = -125029970
# Row 136
$ws.Range("H136").Value = 6262.5757
$ws.Range("I136").Value = 5961.5
$ws.Range("J136").Value = 7065.4443
$ws.Range("K136").Value = 17884.5
$ws.Range("L136").Value = 21196.3329
$ws.Range("M136").Value = -15334.5
$ws.Range("N136").Value = -26296.3329
